$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (44883 = new date) is inserted as row 11,
# pushing the existing rows 11-20 down to 12-21.
$ws.Rows.Item(11).Insert()

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat

$ws.Cells.Item(11, 1).Value  = 6
$ws.Cells.Item(11, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(11, 3).Value  = "Metropolitana"
$ws.Cells.Item(11, 4).Value  = 44883
$ws.Cells.Item(11, 5).Value  = 13
$ws.Cells.Item(11, 6).Value  = 100112010
$ws.Cells.Item(11, 7).Value  = "Achicoria"
$ws.Cells.Item(11, 8).Value  = "Sin especificar"
$ws.Cells.Item(11, 9).Value  = "Primera"
$ws.Cells.Item(11, 10).Value = 180
$ws.Cells.Item(11, 11).Value = 7000
$ws.Cells.Item(11, 12).Value = 8000
$ws.Cells.Item(11, 13).Value = 7500
$ws.Cells.Item(11, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(11, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(11, 16).Value = 469
$ws.Cells.Item(11, 17).Value = 16
$ws.Cells.Item(11, 18).Value = "Hortaliza"
